$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "What's the maximum number of lithology types in an log?"
$ws.Range("B22").Value = "llama3.2:latest"
$ws.Range("C22").Value = "The maximum number of lithology types that can be recorded in a log is 450."
